$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.667659044265747
$ws.Range("B1").Value = 4.846035003662109
$ws.Range("C1").Value = 6.860512733459473
$ws.Range("D1").Value = 6.696332454681396
$ws.Range("E1").Value = 5.300371170043945
